$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 249
$ws.Range("I2").Value = 707
$ws.Range("J2").Value = 2981
$ws.Range("K2").Value = 19
$ws.Range("L2").Value = 761
$ws.Range("M2").Value = 57
$ws.Range("N2").Value = 505
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 9
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 36
$ws.Range("T2").Value = 510
$ws.Range("U2").Value = 32
$ws.Range("V2").Value = 4571
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 4584
$ws.Range("Y2").Value = 6
$ws.Range("Z2").Value = 66
$ws.Range("AA2").Value = 36
